$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern used throughout: cell values that are digit-strings (prices) are
# written via a `="literal"` formula + Copy/PasteSpecial(values-only) round-trip rather
# than a direct .Value assignment, so Excel keeps them as exact text (matching the
# source inlineStr cells) instead of silently parsing look-alike numbers (dropping
# trailing zeros / introducing floating point noise).

# --- Rows 2-36: price/volume updates ---
$ws.Range('D2').Formula = '="73.112.81"'
$ws.Range('D2').Copy()
$ws.Range('D2').PasteSpecial(-4163)
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').Formula = '="4.005.54"'
$ws.Range('D3').Copy()
$ws.Range('D3').PasteSpecial(-4163)
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('D4').Formula = '="0.999"'
$ws.Range('D4').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Formula = '="591.26"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +8.92%  '
$ws.Range('D6').Formula = '="161.19"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +7.83%  '
$ws.Range('D7').Formula = '="0.687"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').Formula = '="0.998"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('E10').Value = '  +1.84%  '
$ws.Range('D11').Formula = '="54.47"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  -3.98%  '
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').Formula = '="11.00"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  +2.95%  '
$ws.Range('D14').Formula = '="4.637.27"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('D15').Formula = '="4.003.53"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  +1.29%  '
$ws.Range('D16').Formula = '="1.27"'
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +8.52%  '
$ws.Range('E17').Value = '  +2.39%  '
$ws.Range('D18').Formula = '="20.47"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').Formula = '="72.861.90"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +2.46%  '
$ws.Range('D21').Formula = '="436.42"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  +2.58%  '
$ws.Range('D22').Formula = '="4.81"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +13.58%  '
$ws.Range('D23').Formula = '="96.46"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  -0.94%  '
$ws.Range('E24').Value = '  -3.60%  '
$ws.Range('D25').Formula = '="4.48"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +19.09%  '
$ws.Range('D26').Formula = '="14.33"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').Formula = '="11.42"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('D28').Formula = '="5.97"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +2.59%  '
$ws.Range('D29').Formula = '="10.51"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('D30').Formula = '="36.51"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('D31').Formula = '="7.93"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  +1.98%  '
$ws.Range('D32').Formula = '="13.78"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +2.78%  '
$ws.Range('E33').Value = '  +1.10%  '
$ws.Range('D34').Formula = '="48.99"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  -5.34%  '
$ws.Range('D35').Formula = '="673.85"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  -3.42%  '
$ws.Range('D36').Formula = '="70.57"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)

# --- Rows 37-40: reordered coins (TheGraph/PEPE swap, ThetaToken/Kaspa swap) ---

# Row 37: PEPE -> TheGraph
$ws.Range('B37').Value = 'TheGraph'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D37').Formula = '="0.442"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  +1.39%  '

# Row 38: TheGraph -> PEPE (price has a subscript-3 char, U+2083, built with UNICHAR
# since PowerShell source is easier to keep ASCII-only; same paste-values trick)
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Formula = '="0.0"&UNICHAR(8323)&"0878"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +6.60%  '

# Row 39: Kaspa -> ThetaToken
$ws.Range('B39').Value = 'ThetaToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D39').Formula = '="3.38"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  -1.56%  '

# Row 40: ThetaToken -> Kaspa
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Formula = '="0.147"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  -2.81%  '

# --- Rows 41-51: price/volume updates ---
$ws.Range('D41').Formula = '="1.00"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  +3.74%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Formula = '="0.0491"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +1.91%  '
$ws.Range('D45').Formula = '="10.86"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +10.71%  '
$ws.Range('E46').Value = '  +1.07%  '
$ws.Range('E47').Value = '  -2.78%  '
$ws.Range('D48').Formula = '="3.39"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('D50').Formula = '="2.839.63"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  +12.62%  '
$ws.Range('D51').Formula = '="3.40"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +4.75%  '
